$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ROADMAP")

# Rows 5 (task #19, "Обработка отчеств и фамилий...") and 9 (task #24,
# "Tree and base context locks") are marked "done": their Priority
# (column E) is cleared and their row formatting is switched to match the
# highlighted "done" style already used by row 12 (task #11).
$doneStyle = $ws.Range("A12:H12")
$doneStyle.Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Range("A9:H9").PasteSpecial(-4122)
$ws.Range("E5").ClearContents()
$ws.Range("E9").ClearContents()

# The sheet keeps its existing sort (Demand descending, then Priority
# ascending); clearing the two Priority values above changes their sort
# position, so re-apply the sort to move them into place.
$sortRange = $ws.Range("A2:H25")
$sortRange.Sort($ws.Range("D2:D25"), 2, $ws.Range("E2:E25"), , 1)

# The re-sort shifts the selection up one row.
$ws.Range("E11").Select()
